$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying data rows (3-16, excluding 2, 8, 14) were re-shuffled: each row's
# Fecha/Volumen/Precio minimo/Precio maximo/Precio promedio ponderado/Origen/Precio $/Kg
# values were replaced with those of another row, per the mapping below.
# Map: target row -> source row (values to copy from, using the *original* data)

$cols = @("D","J","K","L","M","O","P")

# Capture original values (before any changes) for the rows that participate in the reshuffle
$sourceRows = @(3,4,5,6,7,9,10,11,12,13,15,16)
$orig = @{}
foreach ($r in $sourceRows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Range("$c$r").Value()
    }
    $orig[$r] = $rowData
}

# target row -> source row mapping
$mapping = @{
    3  = 11
    4  = 7
    5  = 10
    6  = 12
    7  = 9
    9  = 16
    10 = 13
    11 = 3
    12 = 15
    13 = 5
    15 = 4
    16 = 6
}

foreach ($targetRow in $mapping.Keys) {
    $srcRow = $mapping[$targetRow]
    $vals = $orig[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$targetRow").Value = $vals[$c]
    }
}
